$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row-625 formatting down onto the 13 new rows (626:638) ---
$ws.Range("A625:I625").Copy() | Out-Null
$ws.Range("A626:I638").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Populate the 13 new wellness entries (date 2025-12-02 / serial 45993) ---

# Row 626: Amir Etien
$ws.Range("A626").Value = 45993
$ws.Range("B626").Value = 'Amir Etien'
$ws.Range("C626").Value = 70
$ws.Range("D626").Value = 7
$ws.Range("E626").Value = 6
$ws.Range("F626").Value = 0
$ws.Range("H626").Value = 6

# Row 627: Sofiane Belle
$ws.Range("A627").Value = 45993
$ws.Range("B627").Value = 'Sofiane Belle'
$ws.Range("C627").Value = 70
$ws.Range("D627").Value = 6
$ws.Range("E627").Value = 3
$ws.Range("F627").Value = 0
$ws.Range("H627").Value = 6

# Row 628: Yoan Zouma
$ws.Range("A628").Value = 45993
$ws.Range("B628").Value = 'Yoan Zouma'
$ws.Range("C628").Value = 70
$ws.Range("D628").Value = 7
$ws.Range("E628").Value = 9
$ws.Range("F628").Value = 5
$ws.Range("G628").Value = 'Dos'
$ws.Range("H628").Value = 7

# Row 629: Karim Belmahi
$ws.Range("A629").Value = 45993
$ws.Range("B629").Value = 'Karim Belmahi'
$ws.Range("C629").Value = 70
$ws.Range("D629").Value = 7
$ws.Range("E629").Value = 6
$ws.Range("F629").Value = 0
$ws.Range("H629").Value = 10

# Row 630: Naim Ighbane
$ws.Range("A630").Value = 45993
$ws.Range("B630").Value = 'Naim Ighbane'
$ws.Range("C630").Value = 70
$ws.Range("D630").Value = 6
$ws.Range("E630").Value = 3
$ws.Range("F630").Value = 6
$ws.Range("G630").Value = 'Genou droit'
$ws.Range("H630").Value = 8

# Row 631: Hedi Nasri
$ws.Range("A631").Value = 45993
$ws.Range("B631").Value = 'Hedi Nasri'
$ws.Range("C631").Value = 70
$ws.Range("D631").Value = 7
$ws.Range("E631").Value = 6
$ws.Range("F631").Value = 3
$ws.Range("G631").Value = 'Ischio'
$ws.Range("H631").Value = 7

# Row 632: Ilan Ihaddadene
$ws.Range("A632").Value = 45993
$ws.Range("B632").Value = 'Ilan Ihaddadene'
$ws.Range("C632").Value = 70
$ws.Range("D632").Value = 8
$ws.Range("E632").Value = 8
$ws.Range("F632").Value = 0
$ws.Range("H632").Value = 6

# Row 633: Emmanuel Valey
$ws.Range("A633").Value = 45993
$ws.Range("B633").Value = 'Emmanuel Valey'
$ws.Range("C633").Value = 70
$ws.Range("D633").Value = 8
$ws.Range("E633").Value = 5
$ws.Range("F633").Value = 0
$ws.Range("H633").Value = 7

# Row 634: Karahali Souaré
$ws.Range("A634").Value = 45993
$ws.Range("B634").Value = 'Karahali Souaré'
$ws.Range("C634").Value = 70
$ws.Range("D634").Value = 6
$ws.Range("E634").Value = 6
$ws.Range("F634").Value = 6
$ws.Range("G634").Value = 'Cheville'
$ws.Range("H634").Value = 8

# Row 635: Mattheo Haon
$ws.Range("A635").Value = 45993
$ws.Range("B635").Value = 'Mattheo Haon'
$ws.Range("C635").Value = 70
$ws.Range("D635").Value = 7
$ws.Range("E635").Value = 5
$ws.Range("F635").Value = 0
$ws.Range("H635").Value = 7

# Row 636: Levy Ndoutoume
$ws.Range("A636").Value = 45993
$ws.Range("B636").Value = 'Levy Ndoutoume'
$ws.Range("C636").Value = 70
$ws.Range("D636").Value = 6
$ws.Range("E636").Value = 7
$ws.Range("F636").Value = 1
$ws.Range("G636").Value = 'Ischio'
$ws.Range("H636").Value = 6

# Row 637: Naim Dhib
$ws.Range("A637").Value = 45993
$ws.Range("B637").Value = 'Naim Dhib'
$ws.Range("C637").Value = 70
$ws.Range("D637").Value = 5
$ws.Range("E637").Value = 7
$ws.Range("F637").Value = 0
$ws.Range("H637").Value = 5

# Row 638: Romain Thunet
$ws.Range("A638").Value = 45993
$ws.Range("B638").Value = 'Romain Thunet'
$ws.Range("C638").Value = 70
$ws.Range("D638").Value = 8
$ws.Range("E638").Value = 5
$ws.Range("F638").Value = 3
$ws.Range("G638").Value = 'Ischio '
$ws.Range("H638").Value = 2

# --- Charge column (I): Volume*Intensite formula across the new rows (as one shared group) ---
$ws.Range("I626:I638").Formula = "=C626*D626"

# --- Update the saved view: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 611
$win.ScrollColumn = 1
$ws.Range("D642").Select() | Out-Null

"edit complete"